$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: swap the Emptys/Emptys part counters between the "E2" (empty group)
# and "H2" (numeric_variable group) merged headers.
$ws.Range("E2").Value = 1
$ws.Range("H2").Value = 2

# Row 4: refreshed statistic values (nominal NAs filling repaired)
$ws.Range("B4").Value = 0.6783525101020478
$ws.Range("D4").Value = 0.3567050202040956
$ws.Range("E4").Value = 0.5754920420347929
$ws.Range("G4").Value = 0.1509840840695857
$ws.Range("H4").Value = 0.7796336996336996
$ws.Range("J4").Value = 0.5592673992673991
